# Rotate the contents of rows 14, 15 and 16 on the active sheet:
#   new row 14 <- old row 16
#   new row 15 <- old row 14
#   new row 16 <- old row 15
#
# This is done with Range.Copy so every cell keeps its original value
# *type* (numbers stay numbers; text that merely looks numeric/date-like -
# e.g. the "Antal"/"Starttid" columns - stays text) instead of re-typing
# values through .Value, which would make Excel "smart" re-interpret
# strings such as "10" or "2023-07-28" as numbers/dates.
#
# Range.Copy() here leaves the destination cell untouched whenever the
# source cell is blank (it does not clear it), so every destination row is
# explicitly Clear()-ed right before the rotated data is copied into it,
# guaranteeing stale values never survive the rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original rows being rotated.
$row14 = $ws.Range("A14:AY14")
$row15 = $ws.Range("A15:AY15")
$row16 = $ws.Range("A16:AY16")

# Scratch rows far below the sheet's used range, to stage a copy of each
# row before anything is overwritten.
$stage14 = $ws.Range("A1000:AY1000")
$stage15 = $ws.Range("A1001:AY1001")
$stage16 = $ws.Range("A1002:AY1002")

# 1. Stage the current contents of rows 14-16.
$row14.Copy($stage14)
$row15.Copy($stage15)
$row16.Copy($stage16)

# 2. Wipe the original rows so blank staged cells really end up blank.
$row14.Clear()
$row15.Clear()
$row16.Clear()

# 3. Write the rotated data back: 14<-16, 15<-14, 16<-15.
$stage16.Copy($row14)
$stage14.Copy($row15)
$stage15.Copy($row16)

# 4. Clean up the scratch rows.
$stage14.Clear()
$stage15.Clear()
$stage16.Clear()
